$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AppControl")
$ws2 = $wb.Worksheets.Item("smoke")

# --- AppControl!B25: add email address + hyperlink (mailto) ---
$rngB25 = $ws1.Range("B25")
$rngB25.Value = "stiyyagura@enhops.com"
$ws1.Hyperlinks.Add($rngB25, "mailto:stiyyagura@enhops.com")
# Hyperlinks.Add() re-stamps formatting onto the cell and registers a brand
# new style record; restore the pre-existing "Hyperlink w/ wrap" look so the
# cell keeps referencing the same (already present) cell style.
$rngB25.Style = "Hyperlink"
$rngB25.WrapText = $true

# --- smoke!B18:B25: Run Flag changed from N to Y ---
for ($r = 18; $r -le 25; $r++) {
    $ws2.Range("B$r").Value = "Y"
}

# --- selection / active sheet bookkeeping ---
# Final state: smoke's selection narrows to A26 (single cell) and it is no
# longer the active tab; AppControl becomes the active tab with A26 selected.
$null = $ws2.Range("A26").Select()
$null = $ws1.Activate()
$null = $ws1.Range("A26").Select()
